# Update the dSF column (F) values with repulled data.
# Mapping of row -> new dSF value
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -8
    3  = -3
    5  = -1
    10 = -2
    12 = -2
    13 = 13
    16 = -3
    19 = 0
    21 = 0
    26 = -4
    29 = 1
    32 = -1
    33 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
